# Project Sample Project is saved. TEST Author: admin. Type: SAVE.
#
# Change: cell B11 on the active sheet ("Rules String Hello (Integer hour)")
# is updated from the text "R40" to the text "1". The new value is written
# as text (not a number) so it keeps being stored as a shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")

# Force text storage so that "1" is not auto-converted to a numeric value -
# this mirrors the original cell which held a shared-string value ("R40").
$cell.NumberFormat = "@"
$cell.Value = "1"
